$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 42647.681400462963
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B6").Value = $true
$ws.Range("C6").Value = 9994.36
$ws.Range("D6").Value = 9949.09
$ws.Range("E6").Value = 104.06
$ws.Range("F6").Value = 103.11

$ws.Range("G6").Value = $true
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$ws.Range("H6").Value = -0.91
$ws.Range("I6").Value = $false

$excel.CutCopyMode = $false
